$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 137
$ws.Range("B137").Value = 7499440
$ws.Range("F137").Value = "Olimpia Asuncion"
$ws.Range("G137").Value = "Libertad Asuncion"
$ws.Range("H137").Value = 1
$ws.Range("I137").Value = 3
$ws.Range("K137").Value = 2.7
$ws.Range("M137").Value = 2.3
$ws.Range("N137").Value = 2.8
$ws.Range("O137").Value = 3.25
$ws.Range("P137").Value = 2.3
$ws.Range("Q137").Value = 0.25
$ws.Range("R137").Value = 1.75
$ws.Range("S137").Value = 2.05
$ws.Range("T137").Value = 2.25
$ws.Range("U137").Value = 1.85
$ws.Range("V137").Value = 1.95
$ws.Range("Y137").Value = 1.3
$ws.Range("AA137").Value = 1.05
$ws.Range("AB137").Value = 0.8500000000000001
$ws.Range("AC137").Value = -1

# Row 138
$ws.Range("B138").Value = 7499442
$ws.Range("F138").Value = "Guarani Asuncion"
$ws.Range("G138").Value = "Sportivo Trinidense"
$ws.Range("H138").Value = 0
$ws.Range("I138").Value = 2
$ws.Range("K138").Value = 2
$ws.Range("M138").Value = 3.3
$ws.Range("N138").Value = 1.909
$ws.Range("O138").Value = 3.4
$ws.Range("P138").Value = 3.6
$ws.Range("Q138").Value = -0.5
$ws.Range("R138").Value = 1.975
$ws.Range("S138").Value = 1.825
$ws.Range("T138").Value = 2.5
$ws.Range("U138").Value = 1.9
$ws.Range("V138").Value = 1.9
$ws.Range("Y138").Value = 2.6
$ws.Range("AA138").Value = 0.825
$ws.Range("AB138").Value = -1
$ws.Range("AC138").Value = 0.8999999999999999

# Row 143
$ws.Range("B143").Value = 7493311
$ws.Range("F143").Value = "General Caballero JLM"
$ws.Range("G143").Value = "Olimpia Asuncion"
$ws.Range("H143").Value = 0
$ws.Range("I143").Value = 1
$ws.Range("J143").Value = "A"
$ws.Range("K143").Value = 3.4
$ws.Range("L143").Value = 3.3
$ws.Range("M143").Value = 2
$ws.Range("N143").Value = 3.2
$ws.Range("O143").Value = 3.25
$ws.Range("P143").Value = 2.1
$ws.Range("Q143").Value = 0.25
$ws.Range("R143").Value = 1.95
$ws.Range("S143").Value = 1.85
$ws.Range("T143").Value = 2.25
$ws.Range("U143").Value = 1.775
$ws.Range("V143").Value = 2.025
$ws.Range("W143").Value = -1
$ws.Range("X143").Value = -1
$ws.Range("Y143").Value = 1.1
$ws.Range("Z143").Value = -1
$ws.Range("AA143").Value = 0.8500000000000001
$ws.Range("AB143").Value = -1
$ws.Range("AC143").Value = 1.025

# Row 144
$ws.Range("B144").Value = 7493312
$ws.Range("F144").Value = "Cerro Porteno"
$ws.Range("G144").Value = "Guarani Asuncion"
$ws.Range("H144").Value = 4
$ws.Range("I144").Value = 0
$ws.Range("J144").Value = "H"
$ws.Range("K144").Value = 1.7
$ws.Range("L144").Value = 3.6
$ws.Range("M144").Value = 4.333
$ws.Range("N144").Value = 1.727
$ws.Range("O144").Value = 3.75
$ws.Range("P144").Value = 4.2
$ws.Range("Q144").Value = -0.5
$ws.Range("R144").Value = 1.8
$ws.Range("S144").Value = 2
$ws.Range("T144").Value = 2.75
$ws.Range("U144").Value = 1.875
$ws.Range("V144").Value = 1.925
$ws.Range("W144").Value = 0.7270000000000001
$ws.Range("X144").Value = -1
$ws.Range("Y144").Value = -1
$ws.Range("Z144").Value = 0.8
$ws.Range("AA144").Value = -1
$ws.Range("AB144").Value = 0.875
$ws.Range("AC144").Value = -1

# Row 145
$ws.Range("B145").Value = 7493433
$ws.Range("F145").Value = "Sportivo Luqueno"
$ws.Range("G145").Value = "Nacional Asuncion"
$ws.Range("H145").Value = 1
$ws.Range("I145").Value = 1
$ws.Range("J145").Value = "D"
$ws.Range("K145").Value = 2.75
$ws.Range("L145").Value = 3.2
$ws.Range("M145").Value = 2.4
$ws.Range("N145").Value = 2.75
$ws.Range("O145").Value = 3.1
$ws.Range("P145").Value = 2.45
$ws.Range("Q145").Value = 0.25
$ws.Range("R145").Value = 1.75
$ws.Range("S145").Value = 2.05
$ws.Range("T145").Value = 2.25
$ws.Range("U145").Value = 2
$ws.Range("V145").Value = 1.8
$ws.Range("W145").Value = -1
$ws.Range("X145").Value = 2.1
$ws.Range("Y145").Value = -1
$ws.Range("Z145").Value = 0.375
$ws.Range("AA145").Value = -0.5
$ws.Range("AB145").Value = -0.5
$ws.Range("AC145").Value = 0.4

# Row 210
$ws.Range("R210").Value = 2
$ws.Range("S210").Value = 1.8
$ws.Range("U210").Value = 1.875
$ws.Range("V210").Value = 1.925

# Row 212
$ws.Range("N212").Value = 4.5
$ws.Range("O212").Value = 4.2
$ws.Range("P212").Value = 1.571
$ws.Range("Q212").Value = 1
$ws.Range("U212").Value = 1.85
$ws.Range("V212").Value = 1.95

# Row 213
$ws.Range("R213").Value = 2.025
$ws.Range("S213").Value = 1.775

